# Update "想去人数" (want-to-go count) values in column F across sheets,
# matching the regenerated data output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1838
$ws1.Range("F3").Value = 407
$ws1.Range("F4").Value = 1498
$ws1.Range("F5").Value = 861
$ws1.Range("F8").Value = 13219
$ws1.Range("F9").Value = 13081
$ws1.Range("F13").Value = 550
$ws1.Range("F15").Value = 649
$ws1.Range("F21").Value = 208

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 93
$ws2.Range("F7").Value = 108

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 13

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1838
$ws4.Range("F4").Value = 407
$ws4.Range("F5").Value = 1498
$ws4.Range("F6").Value = 861
$ws4.Range("F8").Value = 93
$ws4.Range("F10").Value = 13219
$ws4.Range("F11").Value = 13081
$ws4.Range("F15").Value = 550
$ws4.Range("F17").Value = 649
$ws4.Range("F27").Value = 13
$ws4.Range("F28").Value = 208
$ws4.Range("F31").Value = 108
